# Rename lab6 -> lab5 references throughout the report (lab05 report).
#
# The document contains several near-duplicate sentences; to avoid touching
# the one sentence that must stay as "lab6-1" ("... для просмотра ..."),
# each Find/Replace below targets the smallest substring that is still
# unique to the run(s) that need to change.

$d = $word.ActiveDocument

# "Открываю файл lfb6-1.asm для редактирования ..." -> lfb5-1
$d.Content.Find.Execute("lfb6-1", $true, $false, $false, $false, $false, `
    $true, 1, $false, "lfb5-1", 2) | Out-Null

# "Создаю объектный файл lab6-1.o, выполняю компоновку ..." -> lab5-1.o
$d.Content.Find.Execute("lab6-1.o", $true, $false, $false, $false, $false, `
    $true, 1, $false, "lab5-1.o", 2) | Out-Null

# "... копирую его из директории Загрузки в директорию lab06 ..." -> lab05
$d.Content.Find.Execute("директорию lab06", $true, $false, $false, $false, `
    $false, $true, 1, $false, "директорию lab05", 2) | Out-Null

# "Копирую файл lab6-1.asm с новым именем lab6-2.asm ..." -> lab5-1 / lab5-2
$d.Content.Find.Execute("Копирую файл lab6-1.asm с новым именем lab6-2.asm", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Копирую файл lab5-1.asm с новым именем lab5-2.asm", 2) | Out-Null

# "Изменяю текст программы lab6-2.asm, чтобы ..." -> lab5-2
$d.Content.Find.Execute("Изменяю текст программы lab6-2.asm", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "Изменяю текст программы lab5-2.asm", 2) | Out-Null

# Both "Создаю объектный файл lab6-2.o, выполняю компоновку ..." sentences
# (one short, one with the extra trailing remark) -> lab5-2.o
$d.Content.Find.Execute("lab6-2.o", $true, $false, $false, $false, $false, `
    $true, 1, $false, "lab5-2.o", 2) | Out-Null

# "В тексте программы lab6-2.asm заменяю sprintLF на sprint ..." -> lab5-2
$d.Content.Find.Execute("В тексте программы lab6-2.asm", $true, $false, `
    $false, $false, $false, $true, 1, $false, `
    "В тексте программы lab5-2.asm", 2) | Out-Null
